# [FIX] budget overview report
#
# Adds seven new filter/criteria rows to the report header block
# (Charge Type, Org, Sector, Subsector, Division, Section, Budget Method)
# right after "Fiscal Year", pushing the existing rows (Functional Area,
# Program Group, Program, Project Group, Project, Run By, Run Date and
# the column-header row) down by 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's very last row (1048576) carries an explicit row-height
# formatting artifact with no cell content. Clear it (and the rest of
# the rows that would otherwise be pushed past the end of the sheet)
# before inserting, so the insert below stays within the sheet bounds.
$ws.Range("A1048570:A1048576").EntireRow.Delete()

# Insert 7 new blank rows starting at row 5, shifting the current
# "Functional Area" row (and everything below it, including the column
# header row) down by 7 rows. Formatting is inherited from row 5.
$ws.Rows.Item(5).Resize(7).Insert()

# Populate the newly inserted rows with the new criteria labels.
$ws.Range("A5").Value = "Charge Type"
$ws.Range("A6").Value = "Org"
$ws.Range("A7").Value = "Sector"
$ws.Range("A8").Value = "Subsector"
$ws.Range("A9").Value = "Division"
$ws.Range("A10").Value = "Section"
$ws.Range("A11").Value = "Budget Method"

# Re-create the blank formatted rows at the very bottom of the sheet
# that were removed above so the sheet ends with the same trailing
# formatting artifact as before the edit.
for ($r = 1048569; $r -le 1048576; $r++) {
    $ws.Rows.Item($r).RowHeight = 12.8
}
